$d = $word.ActiveDocument

# 1) "Épée " + "lourde" -> single run "Épée lourde"
$d.Content.Find.Execute("Épée lourde", $true, $false, $false, $false, $false, $true, 1, $false, "Épée lourde", 2) | Out-Null

# 2) "1,5" + " kg" -> single run "1,5 kg"
$d.Content.Find.Execute("1,5 kg", $true, $false, $false, $false, $false, $true, 1, $false, "1,5 kg", 2) | Out-Null

# 3) "1" + " m" -> single run "1 m"
$d.Content.Find.Execute("1 m", $true, $false, $false, $false, $false, $true, 1, $false, "1 m", 2) | Out-Null

# 4) "2" + "d" + "8+4" -> single run "2d8+4"
$d.Content.Find.Execute("2d8+4", $true, $false, $false, $false, $false, $true, 1, $false, "2d8+4", 2) | Out-Null

# 5) "2" + " PO," + " 7 PA," + " 8 PC" -> single run "2 PO, 7 PA, 8 PC"
$d.Content.Find.Execute("2 PO, 7 PA, 8 PC", $true, $false, $false, $false, $false, $true, 1, $false, "2 PO, 7 PA, 8 PC", 2) | Out-Null

# 6) "Cette épée là est..." -> hyphenate "épée là" to "épée-là" and split the
#    run into three: "Cette " / "épée-là" / " est faites spécialement ..."
$d.Content.Find.Execute("épée là", $true, $false, $false, $false, $false, $true, 1, $false, "épée-là", 2) | Out-Null

$full = $d.Content
$full.Find.Execute("Cette épée-là", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $full.Start

# Force a run split after "Cette " (offset 6) by toggling a formatting
# property on the sub-range (set then restore), which splits the run
# without otherwise altering its properties.
$r1 = $d.Range($paraStart, $paraStart + 6)
$r1.Font.Bold = $true
$r1.Font.Bold = $false

# Force a run split after "épée-là" (offset 13) the same way.
$r2 = $d.Range($paraStart + 6, $paraStart + 13)
$r2.Font.Bold = $true
$r2.Font.Bold = $false
